# Weekly fruit/vegetable price update: two new price records were reported
# for "Terminal Hortofrutícola Agro Chillán" (Tomate, Larga vida, bandeja 20
# kilos) on 2023-08-09 (serial 45147). They are inserted right after the
# existing row 774, pushing every subsequent record down by two rows
# (775->777 ... 826->828), exactly as shown by the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at 775-776; Excel shifts rows 775..826 down to 777..828
# and carries the column D date-style (s="2") onto the new blank rows.
$ws.Rows("775:776").Insert()

# New row 775: Tomate, Larga vida, Primera - $/bandeja 20 kilos
$ws.Cells.Item(775, 1).Value = 7
$ws.Cells.Item(775, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(775, 3).Value = "Ñuble"
$ws.Cells.Item(775, 4).Value = 45147
$ws.Cells.Item(775, 5).Value = 16
$ws.Cells.Item(775, 6).Value = 100112020
$ws.Cells.Item(775, 7).Value = "Tomate"
$ws.Cells.Item(775, 8).Value = "Larga vida"
$ws.Cells.Item(775, 9).Value = "Primera"
$ws.Cells.Item(775, 10).Value = 100
$ws.Cells.Item(775, 11).Value = 18000
$ws.Cells.Item(775, 12).Value = 18000
$ws.Cells.Item(775, 13).Value = 18000
$ws.Cells.Item(775, 14).Value = "$/bandeja 20 kilos"
$ws.Cells.Item(775, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(775, 16).Value = 900
$ws.Cells.Item(775, 17).Value = 20
$ws.Cells.Item(775, 18).Value = "Hortaliza"

# New row 776: Tomate, Larga vida, Segunda - $/bandeja 20 kilos
$ws.Cells.Item(776, 1).Value = 7
$ws.Cells.Item(776, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(776, 3).Value = "Ñuble"
$ws.Cells.Item(776, 4).Value = 45147
$ws.Cells.Item(776, 5).Value = 16
$ws.Cells.Item(776, 6).Value = 100112020
$ws.Cells.Item(776, 7).Value = "Tomate"
$ws.Cells.Item(776, 8).Value = "Larga vida"
$ws.Cells.Item(776, 9).Value = "Segunda"
$ws.Cells.Item(776, 10).Value = 100
$ws.Cells.Item(776, 11).Value = 15000
$ws.Cells.Item(776, 12).Value = 15000
$ws.Cells.Item(776, 13).Value = 15000
$ws.Cells.Item(776, 14).Value = "$/bandeja 20 kilos"
$ws.Cells.Item(776, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(776, 16).Value = 750
$ws.Cells.Item(776, 17).Value = 20
$ws.Cells.Item(776, 18).Value = "Hortaliza"
